# Applies the "New agri dist calculations, backup of script" edit:
#  - Fill in the new "vegheight" (column V) values for the rows that were
#    missing them.
#  - Add two new (essentially blank, date-formatted) trailing rows 32 & 33
#    under the "date" column (E), matching the existing E-column format.
#  - Un-hide / re-size columns D:Q (previously hidden helper columns),
#    leaving column B untouched.
#  - Move the active selection to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "vegheight" (column V) observations for rows 2-30 that did not
#    already have one.
# ---------------------------------------------------------------------
$vegheight = @{
    2  = 1
    3  = 2
    4  = 1.5
    5  = 1
    7  = 1
    8  = 1.2
    9  = 2
    11 = 0.2
    12 = 0.5
    13 = 1.5
    14 = 1.5
    15 = 1
    16 = 2
    18 = 1
    20 = 0.5
    22 = 1.5
    24 = 0.5
    25 = 0.2
    26 = 0.4
    27 = 0.5
    29 = 1
    30 = 0.2
}

foreach ($row in $vegheight.Keys) {
    $ws.Cells.Item($row, 22).Value = $vegheight[$row]
}

# ---------------------------------------------------------------------
# 2. Two new trailing rows carrying the "date" column's format forward
#    (blank entries, e.g. pasted-down formatting from the row above).
# ---------------------------------------------------------------------
$ws.Range("E31").Copy() | Out-Null
$ws.Range("E32").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Un-hide columns D:Q and resize them to fit their (now visible)
#    contents; column B is left at its default width.
# ---------------------------------------------------------------------
$ws.Range("D1:Q1").EntireColumn.Hidden = $false

$colWidths = @{
    1  = 8
    3  = 9.5703125
    4  = 6.5703125
    5  = 9.85546875
    6  = 6
    7  = 7.7109375
    8  = 9.85546875
    9  = 8
    10 = 8
    11 = 5.28515625
    12 = 6.28515625
    13 = 6.28515625
    14 = 4.85546875
    15 = 8.42578125
    16 = 8.28515625
    17 = 4.5703125
    18 = 6.140625
}

foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col] - (5 / 7)
}

# ---------------------------------------------------------------------
# 4. Move the selection to B6 (matches the saved view state).
# ---------------------------------------------------------------------
$ws.Range("B6").Select() | Out-Null

"done"
